$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 11938.95971043
$ws.Range("C2").Value = 11695.2251419348
$ws.Range("E2").Value = 7290.93876967659
$ws.Range("F2").Value = -18.0706703495268

# Row 3
$ws.Range("B3").Value = 12770.186872428
$ws.Range("C3").Value = 12470.892943221
$ws.Range("E3").Value = 8661.11461978616
$ws.Range("F3").Value = 377.339481791966

# Row 4
$ws.Range("B4").Value = 12750.8039240434
$ws.Range("C4").Value = 12482.7728257202
$ws.Range("E4").Value = 8650.39377187479
$ws.Range("F4").Value = 377.387774899792

# Row 5
$ws.Range("B5").Value = 12620.3075662353
$ws.Range("C5").Value = 11602.7241503466
$ws.Range("E5").Value = 8538.47376655087
$ws.Range("F5").Value = 336.055746537395

# Row 6
$ws.Range("B6").Value = 4857.85297688467
$ws.Range("C6").Value = 7742.64575184014
$ws.Range("E6").Value = 7969.22251163519
$ws.Range("F6").Value = 151.500344311472

# Row 7
$ws.Range("B7").Value = 5004.16534219845
$ws.Range("C7").Value = 7961.39718866213
$ws.Range("E7").Value = 8261.28842815431
$ws.Range("F7").Value = 172.784400700685

# Row 8
$ws.Range("C8").Value = 11869.8234422725
$ws.Range("F8").Value = 337.844000043652

# Row 9
$ws.Range("C9").Value = 12377.2444190955
$ws.Range("F9").Value = 358.986540744608

# Row 10
$ws.Range("C10").Value = 12467.6199689492
$ws.Range("F10").Value = 362.752188655181

# Row 11
$ws.Range("C11").Value = 12493.587063291
$ws.Range("F11").Value = 363.834150919421

# Row 12
$ws.Range("C12").Value = 11921.413304804
$ws.Range("F12").Value = 339.99357764913

# Row 13
$ws.Range("C13").Value = 8690.87202853613
$ws.Range("F13").Value = 189.38511676059

# Row 14
$ws.Range("C14").Value = 8225.37311442766
$ws.Range("F14").Value = 265.491839231001

# Row 15
$ws.Range("C15").Value = 11258.6708484081
$ws.Range("F15").Value = 388.36728087083
